# Restore the "From" value of rule R30 (cell C10 on the Rules sheet) from 18 to 1,
# matching the restored revision referenced in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
